$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.333.25'
$ws.Range("E2").Value = '  +12.59%  '
$ws.Range("D3").Value = '1.825.01'
$ws.Range("E3").Value = '  +9.09%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'230.19"
$ws.Range("E5").Value = '  +4.81%  '
$ws.Range("E6").Value = '  +4.22%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = "'31.69"
$ws.Range("E8").Value = '  +6.74%  '
$ws.Range("D9").Value = "'46.97"
$ws.Range("E9").Value = '  +6.38%  '
$ws.Range("E10").Value = '  +7.60%  '
$ws.Range("D11").Value = "'0.0675"
$ws.Range("E11").Value = '  +6.13%  '
$ws.Range("D12").Value = "'0.0930"
$ws.Range("E12").Value = '  +2.89%  '
$ws.Range("D13").Value = '2.086.89'
$ws.Range("E13").Value = '  +9.12%  '
$ws.Range("D14").Value = '1.829.65'
$ws.Range("E14").Value = '  +9.32%  '
$ws.Range("D15").Value = "'0.647"
$ws.Range("E15").Value = '  +5.68%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = "'10.44"
$ws.Range("E16").Value = '  +3.51%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '34.287.93'
$ws.Range("E17").Value = '  +12.34%  '
$ws.Range("D18").Value = "'4.29"
$ws.Range("E18").Value = '  +7.96%  '
$ws.Range("D19").Value = "'70.00"
$ws.Range("E19").Value = '  +5.88%  '
$ws.Range("D20").Value = "'259.24"
$ws.Range("E20").Value = '  +6.97%  '
$ws.Range("D21").Value = '0.0₃0755'
$ws.Range("E21").Value = '  +4.99%  '
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("E23").Value = '  +6.71%  '
$ws.Range("E24").Value = '  +2.52%  '
$ws.Range("E25").Value = '  +3.21%  '
$ws.Range("D26").Value = "'159.05"
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("D27").Value = "'16.68"
$ws.Range("E27").Value = '  +5.38%  '
$ws.Range("D28").Value = "'7.18"
$ws.Range("E28").Value = '  +7.61%  '
$ws.Range("E29").Value = '  +2.59%  '
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("D31").Value = "'3.88"
$ws.Range("E31").Value = '  +12.13%  '
$ws.Range("D32").Value = "'0.0522"
$ws.Range("E32").Value = '  +5.79%  '
$ws.Range("E33").Value = '  +6.04%  '
$ws.Range("D34").Value = "'3.57"
$ws.Range("E34").Value = '  +8.55%  '
$ws.Range("D35").Value = '1.555.08'
$ws.Range("E35").Value = '  +3.67%  '
$ws.Range("D36").Value = "'1.81"
$ws.Range("E36").Value = '  +3.27%  '
$ws.Range("E37").Value = '  +7.02%  '
$ws.Range("D38").Value = "'85.70"
$ws.Range("E38").Value = '  +1.80%  '
$ws.Range("D39").Value = "'0.632"
$ws.Range("E39").Value = '  +6.24%  '
$ws.Range("E40").Value = '  +5.76%  '
$ws.Range("D41").Value = "'2.82"
$ws.Range("E41").Value = '  +5.47%  '
$ws.Range("D42").Value = "'0.921"
$ws.Range("E42").Value = '  +9.95%  '
$ws.Range("D43").Value = "'2.33"
$ws.Range("E43").Value = '  +1.83%  '
$ws.Range("D44").Value = "'2.18"
$ws.Range("E44").Value = '  +10.78%  '
$ws.Range("D45").Value = "'0.0526"
$ws.Range("E45").Value = '  +5.76%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = "'1.07"
$ws.Range("E46").Value = '  +4.80%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.987.47'
$ws.Range("E47").Value = '  +10.13%  '
$ws.Range("B48").Value = 'MinaProtocolToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range("D48").Value = "'1.07"
$ws.Range("E48").Value = '  +160.72%  '
$ws.Range("D49").Value = "'12.51"
$ws.Range("E49").Value = '  +29.77%  '
$ws.Range("D50").Value = "'5.76"
$ws.Range("E50").Value = '  +4.26%  '
$ws.Range("D51").Value = "'53.15"
$ws.Range("E51").Value = '  +3.94%  '
